$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column, new cell text. The source data keeps these
# "Price"/"Volume(1h)" figures as plain text (as scraped), so any value that
# merely looks like a plain decimal number is written with a Text number
# format first -- otherwise Excel would silently reinterpret it as a number
# (e.g. "222.88" -> 222.88) instead of keeping the original string "222.88".
$updates = @(
    ,@(2, 4, '34.135.99')
    ,@(2, 5, '  -1.58%  ')
    ,@(3, 4, '1.797.65')
    ,@(3, 5, '  +0.54%  ')
    ,@(4, 5, '  -0.18%  ')
    ,@(5, 4, '222.88')
    ,@(5, 5, '  +0.07%  ')
    ,@(6, 4, '0.551')
    ,@(6, 5, '  -0.52%  ')
    ,@(7, 5, '  -0.11%  ')
    ,@(8, 4, '32.31')
    ,@(8, 5, '  -0.70%  ')
    ,@(9, 4, '0.285')
    ,@(9, 5, '  +2.09%  ')
    ,@(10, 4, '0.0718')
    ,@(10, 5, '  +5.20%  ')
    ,@(11, 5, '  -1.44%  ')
    ,@(12, 4, '2.055.95')
    ,@(12, 5, '  +0.48%  ')
    ,@(13, 4, '1.791.50')
    ,@(13, 5, '  +0.34%  ')
    ,@(14, 4, '10.72')
    ,@(14, 5, '  -2.60%  ')
    ,@(15, 4, '0.631')
    ,@(15, 5, '  +0.33%  ')
    ,@(16, 4, '34.098.81')
    ,@(16, 5, '  -1.73%  ')
    ,@(17, 5, '  -1.48%  ')
    ,@(18, 4, '68.09')
    ,@(18, 5, '  -0.41%  ')
    ,@(19, 4, '246.44')
    ,@(19, 5, '  -2.41%  ')
    ,@(20, 4, '0.0₃0785')
    ,@(20, 5, '  +0.79%  ')
    ,@(21, 5, '  +0.02%  ')
    ,@(22, 4, '10.78')
    ,@(22, 5, '  +3.08%  ')
    ,@(23, 4, '4.10')
    ,@(23, 5, '  -1.81%  ')
    ,@(24, 5, '  -1.23%  ')
    ,@(26, 4, '16.49')
    ,@(26, 5, '  +0.97%  ')
    ,@(27, 5, '  +0.32%  ')
    ,@(28, 5, '  -1.47%  ')
    ,@(29, 5, '  -0.20%  ')
    ,@(30, 4, '0.0519')
    ,@(30, 5, '  +0.88%  ')
    ,@(31, 2, 'PancakeSwap')
    ,@(31, 3, 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake')
    ,@(31, 4, '1.21')
    ,@(31, 5, '  +1.95%  ')
    ,@(32, 2, 'Filecoin')
    ,@(32, 3, 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil')
    ,@(32, 4, '3.71')
    ,@(32, 5, '  -0.61%  ')
    ,@(33, 5, '  -1.31%  ')
    ,@(34, 4, '1.85')
    ,@(34, 5, '  -0.52%  ')
    ,@(35, 4, '1.415.58')
    ,@(35, 5, '  -1.22%  ')
    ,@(36, 4, '0.645')
    ,@(36, 5, '  +2.50%  ')
    ,@(37, 5, '  +0.43%  ')
    ,@(38, 5, '  -1.27%  ')
    ,@(39, 4, '0.945')
    ,@(39, 5, '  +4.98%  ')
    ,@(40, 4, '80.27')
    ,@(40, 5, '  -2.85%  ')
    ,@(41, 5, '  -2.65%  ')
    ,@(42, 5, '  -0.40%  ')
    ,@(43, 4, '2.14')
    ,@(43, 5, '  +4.07%  ')
    ,@(44, 5, '  +0.23%  ')
    ,@(45, 4, '0.0496')
    ,@(45, 5, '  -1.30%  ')
    ,@(46, 4, '1.954.80')
    ,@(46, 5, '  +0.91%  ')
    ,@(47, 5, '  -2.61%  ')
    ,@(48, 4, '106.27')
    ,@(48, 5, '  +1.86%  ')
    ,@(49, 5, '  -0.11%  ')
    ,@(50, 4, '11.90')
    ,@(50, 5, '  -0.25%  ')
    ,@(51, 4, '0.0₆0122')
    ,@(51, 5, '  +0.47%  ')
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $cell = $ws.Cells.Item($row, $col)
    if ($val -match '^[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $val
}
